$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '61.933.68'
$ws.Range('E2').Value = '  -1.19%  '
$ws.Range('D3').Value = '3.406.85'
$ws.Range('E3').Value = '  -0.73%  '
$ws.Range('E4').Value = '  +0.00%  '
$ws.Range('D5').Value = '''409.14'
$ws.Range('E5').Value = '  +0.52%  '
$ws.Range('D6').Value = '''129.04'
$ws.Range('E6').Value = '  -1.06%  '
$ws.Range('E7').Value = '  +6.24%  '
$ws.Range('D9').Value = '''0.732'
$ws.Range('E9').Value = '  +5.49%  '
$ws.Range('E10').Value = '  +2.91%  '
$ws.Range('D11').Value = '''42.78'
$ws.Range('E11').Value = '  +1.87%  '
$ws.Range('D12').Value = '''0.0000222'
$ws.Range('E12').Value = '  +40.46%  '
$ws.Range('D13').Value = '''9.30'
$ws.Range('E13').Value = '  +10.22%  '
$ws.Range('E14').Value = '  -0.43%  '
$ws.Range('D15').Value = '3.950.40'
$ws.Range('E15').Value = '  -0.85%  '
$ws.Range('E16').Value = '  +6.91%  '
$ws.Range('D17').Value = '3.410.15'
$ws.Range('E17').Value = '  -1.04%  '
$ws.Range('D18').Value = '''12.53'
$ws.Range('E18').Value = '  +8.94%  '
$ws.Range('E19').Value = '  +6.94%  '
$ws.Range('D20').Value = '61.931.09'
$ws.Range('E20').Value = '  -1.22%  '
$ws.Range('D21').Value = '''448.52'
$ws.Range('E21').Value = '  +42.44%  '
$ws.Range('D22').Value = '''91.98'
$ws.Range('E22').Value = '  +8.74%  '
$ws.Range('D23').Value = '''3.21'
$ws.Range('E23').Value = '  +0.89%  '
$ws.Range('D24').Value = '''13.21'
$ws.Range('E24').Value = '  +2.98%  '
$ws.Range('D25').Value = '''3.29'
$ws.Range('E25').Value = '  +3.70%  '
$ws.Range('E26').Value = '  +14.52%  '
$ws.Range('D27').Value = '''33.08'
$ws.Range('E27').Value = '  +11.24%  '
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('D29').Value = '''7.77'
$ws.Range('E29').Value = '  -0.84%  '
$ws.Range('E30').Value = '  +0.23%  '
$ws.Range('D31').Value = '''11.98'
$ws.Range('E31').Value = '  +5.11%  '
$ws.Range('E32').Value = '  -1.83%  '
$ws.Range('E33').Value = '  -0.14%  '
$ws.Range('D34').Value = '''42.74'
$ws.Range('E34').Value = '  -3.95%  '
$ws.Range('E35').Value = '  +0.02%  '
$ws.Range('D36').Value = '''0.0504'
$ws.Range('E36').Value = '  +4.22%  '
$ws.Range('D37').Value = '''53.77'
$ws.Range('E37').Value = '  +3.65%  '
$ws.Range('D38').Value = '''0.998'
$ws.Range('E38').Value = '  -0.09%  '
$ws.Range('E39').Value = '  +1.74%  '
$ws.Range('E40').Value = '  +7.74%  '
$ws.Range('B41').Value = 'TheGraph'
$ws.Range('C41').Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range('D41').Value = '''0.321'
$ws.Range('E41').Value = '  -0.52%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''2.94'
$ws.Range('E42').Value = '  -0.68%  '
$ws.Range('D43').Value = '''143.28'
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '''4.27'
$ws.Range('E44').Value = '  +9.17%  '
$ws.Range('D45').Value = '''2.56'
$ws.Range('E45').Value = '  +15.51%  '
$ws.Range('E46').Value = '  +0.88%  '
$ws.Range('E47').Value = '  -1.54%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '''22.29'
$ws.Range('E48').Value = '  +4.90%  '
$ws.Range('B49').Value = 'Cronos'
$ws.Range('C49').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D49').Value = '''0.146'
$ws.Range('E49').Value = '  +20.66%  '
$ws.Range('D50').Value = '''2.14'
$ws.Range('E50').Value = '  +8.18%  '
$ws.Range('D51').Value = '''1.91'
$ws.Range('E51').Value = '  +13.08%  '
